# "Modificación de permisos terminado"
# Adds the missing "Crear Partida" / "Modificar Partida" permission rows to
# the Patentes table (matching the pattern already used for the other ABM
# entities) and colors the "ABM Partida" row in the Familias table to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New permission rows for "ABM Partida" -------------------------------
$ws.Range("M15").Value = 9
$ws.Range("N15").Value = "Crear Partida"
$ws.Range("M16").Value = 10
$ws.Range("N16").Value = "Modificar Partida"

# Match the numbering column's existing look (vertical-centered, wrapped) -
# same formatting already used by M7:M14.
$ws.Range("M15:M16").VerticalAlignment = -4108
$ws.Range("M15:M16").WrapText = $true

# Color the new rows with the same highlight color used for "ABM Partida"
# (a light blue-gray — Theme color "Blue-Gray, Text 2, Lighter 60%").
$partidaColor = 13285805
$ws.Range("N15:N16").Interior.Color = $partidaColor

# Highlight the existing "ABM Partida" cell to match, completing the color
# coding already present for the other ABM entries.
$ws.Range("I9").Interior.Color = $partidaColor

# Scroll the view down a bit so the newly added rows are visible.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 4
